$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.38"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "25.40"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.130"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05588"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.511"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8170"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8411"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1344"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06949"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03159"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02846"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09378"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001515"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0005953"
$ws.Range("E16").Value = "15OneONE"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006090"
$ws.Range("E17").Value = "16TigerCashTCH"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.501"
$ws.Range("E18").Value = "17LEOLEO"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.063"
$ws.Range("E19").Value = "18BTSETokenBTSE"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3179"
$ws.Range("E20").Value = "19BitpandaEcosystemTokenBEST"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.746"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04733"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1376"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001247"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004647"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009704"
$ws.Range("E27").Value = "26NitroExNTXBestin24h"
$ws.Range("E28").Value = "27UpBotsUBXT"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03665"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1365"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002625"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003378"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008284"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005298"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002124"
